{"js": "// Update the date line and the 25 division \"problem, answer\" cells.\n// Each old text is unique within the document, so plain search + replace\n// (matchCase, no wildcards) is safe and order-independent.\nconst replacements = [\n  [\"2024-02-09 Friday\", \"2024-02-10 Saturday\"],\n  [\"77\u00f74=19, 1\", \"50\u00f73=16, 2\"],\n  [\"80\u00f73=26, 2\", \"53\u00f72=26, 1\"],\n  [\"59\u00f74=14, 3\", \"26\u00f74=6, 2\"],\n  [\"62\u00f73=20, 2\", \"56\u00f74=14, 0\"],\n  [\"56\u00f72=28, 0\", \"11\u00f79=1, 2\"],\n  [\"58\u00f73=19, 1\", \"84\u00f72=42, 0\"],\n  [\"83\u00f72=41, 1\", \"65\u00f73=21, 2\"],\n  [\"45\u00f72=22, 1\", \"92\u00f79=10, 2\"],\n  [\"70\u00f76=11, 4\", \"28\u00f72=14, 0\"],\n  [\"82\u00f76=13, 4\", \"91\u00f73=30, 1\"],\n  [\"21\u00f74=5, 1\", \"86\u00f73=28, 2\"],\n  [\"61\u00f78=7, 5\", \"53\u00f79=5, 8\"],\n  [\"59\u00f77=8, 3\", \"55\u00f74=13, 3\"],\n  [\"27\u00f75=5, 2\", \"44\u00f79=4, 8\"],\n  [\"51\u00f74=12, 3\", \"39\u00f79=4, 3\"],\n  [\"39\u00f75=7, 4\", \"87\u00f79=9, 6\"],\n  [\"93\u00f76=15, 3\", \"43\u00f73=14, 1\"],\n  [\"14\u00f75=2, 4\", \"98\u00f76=16, 2\"],\n  [\"20\u00f78=2, 4\", \"29\u00f74=7, 1\"],\n  [\"88\u00f78=11, 0\", \"88\u00f77=12, 4\"],\n  [\"71\u00f76=11, 5\", \"50\u00f77=7, 1\"],\n  [\"69\u00f74=17, 1\", \"44\u00f73=14, 2\"],\n  [\"76\u00f72=38, 0\", \"60\u00f78=7, 4\"],\n  [\"77\u00f77=11, 0\", \"71\u00f79=7, 8\"],\n  [\"56\u00f78=7, 0\", \"40\u00f78=5, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division \"problem, answer\" cells.\n# Each old text is unique within the document, so Find/Replace with\n# MatchCase + Replace:=wdReplaceAll (2) is safe and order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-09 Friday\", \"2024-02-10 Saturday\"),\n    @(\"77\u00f74=19, 1\", \"50\u00f73=16, 2\"),\n    @(\"80\u00f73=26, 2\", \"53\u00f72=26, 1\"),\n    @(\"59\u00f74=14, 3\", \"26\u00f74=6, 2\"),\n    @(\"62\u00f73=20, 2\", \"56\u00f74=14, 0\"),\n    @(\"56\u00f72=28, 0\", \"11\u00f79=1, 2\"),\n    @(\"58\u00f73=19, 1\", \"84\u00f72=42, 0\"),\n    @(\"83\u00f72=41, 1\", \"65\u00f73=21, 2\"),\n    @(\"45\u00f72=22, 1\", \"92\u00f79=10, 2\"),\n    @(\"70\u00f76=11, 4\", \"28\u00f72=14, 0\"),\n    @(\"82\u00f76=13, 4\", \"91\u00f73=30, 1\"),\n    @(\"21\u00f74=5, 1\", \"86\u00f73=28, 2\"),\n    @(\"61\u00f78=7, 5\", \"53\u00f79=5, 8\"),\n    @(\"59\u00f77=8, 3\", \"55\u00f74=13, 3\"),\n    @(\"27\u00f75=5, 2\", \"44\u00f79=4, 8\"),\n    @(\"51\u00f74=12, 3\", \"39\u00f79=4, 3\"),\n    @(\"39\u00f75=7, 4\", \"87\u00f79=9, 6\"),\n    @(\"93\u00f76=15, 3\", \"43\u00f73=14, 1\"),\n    @(\"14\u00f75=2, 4\", \"98\u00f76=16, 2\"),\n    @(\"20\u00f78=2, 4\", \"29\u00f74=7, 1\"),\n    @(\"88\u00f78=11, 0\", \"88\u00f77=12, 4\"),\n    @(\"71\u00f76=11, 5\", \"50\u00f77=7, 1\"),\n    @(\"69\u00f74=17, 1\", \"44\u00f73=14, 2\"),\n    @(\"76\u00f72=38, 0\", \"60\u00f78=7, 4\"),\n    @(\"77\u00f77=11, 0\", \"71\u00f79=7, 8\"),\n    @(\"56\u00f78=7, 0\", \"40\u00f78=5, 0\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
